# "Add existing VRES and BESS"
# Update the "Power Storage" sheet: set already-existing BESS units (ExisUnits,
# column E) for the new nodes, and reduce the remaining investable capacity
# (MaxInvest, column S) accordingly for rows 7-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ExisUnits (column E)
$ws.Range("E7").Value = 16
$ws.Range("E9").Value = 40

# MaxInvest (column S)
$ws.Range("S7").Value = 8
$ws.Range("S8").Value = 8
$ws.Range("S9").Value = 8
$ws.Range("S10").Value = 8
$ws.Range("S11").Value = 8

# Move the selection in the frozen (bottom-left) pane to S12, matching the
# cell the author last interacted with after entering the new values.
$ws.Activate()
$ws.Range("S12").Select()
